$wb = $excel.ActiveWorkbook

# Update "想去人数" (F) counts on both the "展览" and "全部类型" sheets,
# as regenerated at commit 456a3b4.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 474
    $ws.Range("F3").Value = 54
}
